$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = "0-Input/PLAST-60/derog cfep0033-e.pdf"
$ws.Range("I3").Value = "0-Input/PLAST-60/cmet0005 derog.pdf"
$ws.Range("I4").Value = "0-Input/PLAST-60/cfep0032 derog.pdf"
$ws.Range("I5").Value = "0-Input/PLAST-60/cfep0065 2022 derog.pdf"
$ws.Range("I6").Value = "0-Input/PLAST-60/derogation ccom0029 32.pdf"
$ws.Range("I7").Value = "0-Input/PLAST-60/derogation ccom0029 32.pdf"
$ws.Range("I8").Value = "0-Input/PLAST-60/derog enceinte.pdf"
$ws.Range("I9").Value = "0-Input/PLAST-60/derog balance 3kg kern.pdf"
$ws.Range("I10").Value = "0-Input/PLAST-60/derog balance 3kg kern.pdf"
$ws.Range("I11").Value = "0-Input/PLAST-60/derog balance 3kg kern.pdf"
$ws.Range("I12").Value = "0-Input/PLAST-60/derog module PEAK.pdf"
$ws.Range("I13").Value = "0-Input/PLAST-60/scr-058 derog.pdf"
$ws.Range("I14").Value = "0-Input/PLAST-60/pv cycling derog.pdf"
$ws.Range("I15").Value = "0-Input/PLAST-60/cemi0016 reformé.pdf"
$ws.Range("I16").Value = "0-Input/PLAST-60/retour cofrac.pdf"
$ws.Range("I17").Value = "0-Input/PLAST-60/mail.pdf"
$ws.Range("I18").Value = "0-Input/PLAST-60/constructeur.pdf"
$ws.Range("I19").Value = "0-Input/PLAST-60/retour banc pendule.pdf"
$ws.Range("I20").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo.pdf"
$ws.Range("I21").Value = "0-Input/PLAST-60/mail.pdf"
$ws.Range("I22").Value = "0-Input/PLAST-60/retour cofrac.pdf"
$ws.Range("I23").Value = "0-Input/PLAST-60/Microsoft Outlook - MémoCMET0044.pdf"
$ws.Range("I24").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo banc lavalll.pdf"
$ws.Range("I25").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo banc lavalll.pdf"
$ws.Range("I26").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo cscr0471.pdf"
$ws.Range("I27").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémocivc0058.pdf"
$ws.Range("I28").Value = "0-Input/PLAST-60/CEPE0014 ET CFEP0049.pdf"
$ws.Range("I29").Value = "0-Input/PLAST-60/CEPE0014 ET CFEP0049.pdf"
$ws.Range("I30").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo THERMOCOUPLE scr.pdf"
$ws.Range("I31").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo cmag0242.pdf"
$ws.Range("I32").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo ni.pdf"
$ws.Range("I33").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo ni.pdf"
$ws.Range("I34").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo ni.pdf"
$ws.Range("I35").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo ni.pdf"
$ws.Range("I36").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo ccom0034.pdf"
$ws.Range("I37").Value = "0-Input/PLAST-60/Fiche renseignement PO Ajout Equipement - Isopar Station.doc"
$ws.Range("I38").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo  FID.pdf"
$ws.Range("I39").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo  FID.pdf"
$ws.Range("I40").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo  FID.pdf"
$ws.Range("I41").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo clab0221.pdf"
$ws.Range("I42").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo  FID.pdf"
$ws.Range("I43").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo  FID.pdf"
$ws.Range("I44").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo alimentation.pdf"
$ws.Range("I45").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo alimentation.pdf"
$ws.Range("I46").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo alimentation.pdf"
$ws.Range("I47").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo alimentation.pdf"
$ws.Range("I48").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo alimentation.pdf"
$ws.Range("I49").Value = "0-Input/PLAST-60/derogation moyen d'approche.pdf"
$ws.Range("I50").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo lvdt derogation.pdf"
$ws.Range("I51").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo lvdt derogation.pdf"
$ws.Range("I52").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo lvdt derogation.pdf"
$ws.Range("I53").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo lvdt derogation.pdf"
$ws.Range("I54").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo lvdt derogation.pdf"
$ws.Range("I55").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo lvdt derogation.pdf"
$ws.Range("I56").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo lvdt derogation.pdf"
$ws.Range("I57").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo lvdt derogation.pdf"
$ws.Range("I58").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo lvdt derogation.pdf"
$ws.Range("I59").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo lvdt derogation.pdf"
$ws.Range("I60").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo lvdt derogation.pdf"
$ws.Range("I61").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo lvdt derogation.pdf"
$ws.Range("I62").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo lvdt derogation.pdf"
$ws.Range("I63").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo lvdt derogation.pdf"
$ws.Range("I64").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo lvdt derogation.pdf"
$ws.Range("I65").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo lvdt derogation.pdf"
$ws.Range("I66").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo lvdt derogation.pdf"
$ws.Range("I67").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo lvdt derogation.pdf"
$ws.Range("I68").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo lvdt derogation.pdf"
$ws.Range("I69").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo lvdt derogation.pdf"
$ws.Range("I70").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo lvdt derogation.pdf"
$ws.Range("I71").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo lvdt derogation.pdf"
$ws.Range("I72").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo lvdt derogation.pdf"
$ws.Range("I73").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo lvdt derogation.pdf"
$ws.Range("I74").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo lvdt derogation.pdf"
$ws.Range("I75").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo lvdt derogation.pdf"
$ws.Range("I76").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo lvdt derogation.pdf"
$ws.Range("I77").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo lvdt derogation.pdf"
$ws.Range("I78").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo lvdt derogation.pdf"
$ws.Range("I79").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo lvdt derogation.pdf"
$ws.Range("I80").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo lvdt derogation.pdf"
$ws.Range("I81").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo lvdt derogation.pdf"
$ws.Range("I82").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo lvdt derogation.pdf"
$ws.Range("I83").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo lvdt derogation.pdf"
$ws.Range("I84").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo lvdt derogation.pdf"
$ws.Range("I85").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo lvdt derogation.pdf"
$ws.Range("I86").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo lvdt derogation.pdf"
$ws.Range("I87").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo lvdt derogation.pdf"
$ws.Range("I88").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo lvdt derogation.pdf"
$ws.Range("I89").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo lvdt derogation.pdf"
$ws.Range("I90").Value = "0-Input/PLAST-60/Microsoft Outlook - Mémo lvdt derogation.pdf"
$ws.Range("I91").Value = "0-Input/PLAST-60/Microsoft Outlook - cmag0242.pdf"
$ws.Range("I92").Value = "0-Input/PLAST-60/Microsoft Outlook - PT 100 station.pdf"
$ws.Range("I93").Value = "0-Input/PLAST-60/Microsoft Outlook - cepe0008.pdf"
$ws.Range("I94").Value = "0-Input/PLAST-60/creation instrument LVDT atex.pdf"
$ws.Range("I95").Value = "0-Input/PLAST-60/creation instrument LVDT atex.pdf"
$ws.Range("I96").Value = "0-Input/PLAST-60/Microsoft Outlook - alim.pdf"
$ws.Range("I97").Value = "0-Input/PLAST-60/Microsoft Outlook - alim.pdf"
$ws.Range("I98").Value = "0-Input/PLAST-60/Microsoft Outlook - alim.pdf"
$ws.Range("I99").Value = "0-Input/PLAST-60/Microsoft Outlook - alim.pdf"
$ws.Range("I100").Value = "0-Input/PLAST-60/Microsoft Outlook - alim.pdf"
$ws.Range("I101").Value = "0-Input/PLAST-60/Microsoft Outlook - alim.pdf"
$ws.Range("I102").Value = "0-Input/PLAST-60/Microsoft Outlook - alim.pdf"
$ws.Range("I103").Value = "0-Input/PLAST-60/Microsoft Outlook - alim.pdf"
$ws.Range("I104").Value = "0-Input/PLAST-60/Microsoft Outlook - alim.pdf"
$ws.Range("I105").Value = "0-Input/PLAST-60/Microsoft Outlook - ccal0080.pdf"
$ws.Range("I106").Value = "0-Input/PLAST-60/derogation FEP-024 (CFEP0201).docx"
$ws.Range("I107").Value = "0-Input/PLAST-60/derogation MET-002.docx"
$ws.Range("I108").Value = "0-Input/PLAST-60/EMT 5 à 6%.msg"
$ws.Range("I109").Value = "0-Input/PLAST-60/derogation sled.pdf"
$ws.Range("I110").Value = "0-Input/PLAST-60/derogation sled.pdf"
Write-Output "Updated 109 validation proof paths"
